$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C entirely (the old "max" column of all-1 values).
# This shifts column D ("prediction") into C and column E ("rejection-f") into D.
$ws.Range("C1:C9").Delete()

# Update column B values (previously all 1) with the actual computed values.
$ws.Range("B2").Value = 22.3185113181394
$ws.Range("B3").Value = 19.41846384564789
$ws.Range("B4").Value = 9.479033267193529
$ws.Range("B5").Value = 2.884868659424686
$ws.Range("B6").Value = 7.269836897346863
$ws.Range("B7").Value = 18.42405118119634
$ws.Range("B8").Value = 20.64837580530254
$ws.Range("B9").Value = 27.7815495458109
